# Update cryptocurrency price/volume figures per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are free-form text (e.g. "41.910.97"); force text
# formatting before/after the write so Excel does not reinterpret values
# that happen to look numeric (e.g. "1.00" -> 1) while leaving the cell
# format back at General, matching the original workbook styling.
$priceCells = @(
    @{Cell='D2'; Value='41.910.97'},
    @{Cell='D3'; Value='2.249.81'},
    @{Cell='D5'; Value='229.75'},
    @{Cell='D7'; Value='61.31'},
    @{Cell='D10'; Value='58.65'},
    @{Cell='D11'; Value='0.0882'},
    @{Cell='D13'; Value='2.582.26'},
    @{Cell='D14'; Value='15.86'},
    @{Cell='D15'; Value='21.73'},
    @{Cell='D16'; Value='0.805'},
    @{Cell='D17'; Value='5.61'},
    @{Cell='D18'; Value='2.249.00'},
    @{Cell='D19'; Value='41.815.04'},
    @{Cell='D23'; Value='250.24'},
    @{Cell='D26'; Value='2.38'},
    @{Cell='D27'; Value='9.64'},
    @{Cell='D29'; Value='168.09'},
    @{Cell='D30'; Value='20.19'},
    @{Cell='D32'; Value='2.80'},
    @{Cell='D34'; Value='5.17'},
    @{Cell='D38'; Value='6.72'},
    @{Cell='D42'; Value='5.00'},
    @{Cell='D44'; Value='8.75'},
    @{Cell='D45'; Value='100.69'},
    @{Cell='D46'; Value='0.0981'},
    @{Cell='D47'; Value='1.488.18'},
    @{Cell='D49'; Value='16.56'}
)
foreach ($item in $priceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.NumberFormat = "General"
}

# Volume(1h) cells (column E) already contain "%"-suffixed text, so a plain
# Value assignment keeps them as text without any coercion to a number.
$volumeCells = @(
    @{Cell='E2'; Value='  +5.55%  '},
    @{Cell='E3'; Value='  +3.97%  '},
    @{Cell='E4'; Value='  +0.13%  '},
    @{Cell='E5'; Value='  +0.75%  '},
    @{Cell='E6'; Value='  -0.52%  '},
    @{Cell='E7'; Value='  -3.55%  '},
    @{Cell='E8'; Value='  +0.09%  '},
    @{Cell='E9'; Value='  +3.19%  '},
    @{Cell='E10'; Value='  +0.95%  '},
    @{Cell='E11'; Value='  +3.74%  '},
    @{Cell='E12'; Value='  +0.35%  '},
    @{Cell='E13'; Value='  +3.94%  '},
    @{Cell='E14'; Value='  -1.27%  '},
    @{Cell='E15'; Value='  -1.34%  '},
    @{Cell='E16'; Value='  -0.58%  '},
    @{Cell='E17'; Value='  +2.09%  '},
    @{Cell='E18'; Value='  +3.94%  '},
    @{Cell='E19'; Value='  +5.56%  '},
    @{Cell='E21'; Value='  +0.33%  '},
    @{Cell='E22'; Value='  +6.03%  '},
    @{Cell='E23'; Value='  +9.70%  '},
    @{Cell='E24'; Value='  +0.06%  '},
    @{Cell='E25'; Value='  +2.39%  '},
    @{Cell='E26'; Value='  -0.79%  '},
    @{Cell='E27'; Value='  +1.52%  '},
    @{Cell='E28'; Value='  +2.56%  '},
    @{Cell='E29'; Value='  -2.52%  '},
    @{Cell='E30'; Value='  +1.97%  '},
    @{Cell='E31'; Value='  +2.09%  '},
    @{Cell='E32'; Value='  +4.20%  '},
    @{Cell='E33'; Value='  +0.55%  '},
    @{Cell='E34'; Value='  +10.15%  '},
    @{Cell='E35'; Value='  +1.88%  '},
    @{Cell='E36'; Value='  +1.49%  '},
    @{Cell='E37'; Value='  +4.59%  '},
    @{Cell='E38'; Value='  -3.53%  '},
    @{Cell='E39'; Value='  -0.41%  '},
    @{Cell='E40'; Value='  +33.41%  '},
    @{Cell='E41'; Value='  +0.13%  '},
    @{Cell='E42'; Value='  +3.68%  '},
    @{Cell='E43'; Value='  +4.95%  '},
    @{Cell='E44'; Value='  +12.93%  '},
    @{Cell='E45'; Value='  -1.39%  '},
    @{Cell='E46'; Value='  +6.02%  '},
    @{Cell='E47'; Value='  -1.58%  '},
    @{Cell='E48'; Value='  -1.65%  '},
    @{Cell='E49'; Value='  -6.05%  '},
    @{Cell='E50'; Value='  +0.12%  '},
    @{Cell='E51'; Value='  -0.51%  '}
)
foreach ($item in $volumeCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
